$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62,8).Value = 8395.333000000001
$ws.Cells.Item(62,9).Value = 8395.333000000001
$ws.Cells.Item(62,10).Value = 0
$ws.Cells.Item(62,11).Value = 8395.333000000001
$ws.Cells.Item(62,12).Value = 0
$ws.Cells.Item(62,13).ClearContents()
$ws.Cells.Item(62,14).Value = -7771.333000000001
$ws.Cells.Item(65,8).Value = 8395.333000000001
$ws.Cells.Item(65,9).Value = 8395.333000000001
$ws.Cells.Item(65,10).Value = 0
$ws.Cells.Item(65,11).Value = 41976.665
$ws.Cells.Item(65,12).Value = 0
$ws.Cells.Item(65,13).ClearContents()
$ws.Cells.Item(65,14).Value = -38856.665
$ws.Cells.Item(131,8).Value = 4282.5
$ws.Cells.Item(131,10).Value = 4180
$ws.Cells.Item(131,12).Value = 12540
$ws.Cells.Item(131,14).Value = -22620
$ws.Cells.Item(137,8).Value = 1080.3334
$ws.Cells.Item(137,9).Value = 1080.3334
$ws.Cells.Item(137,11).Value = 3241.0002
$ws.Cells.Item(137,13).Value = -691.0001999999999
$ws.Cells.Item(138,8).Value = 2336.516
$ws.Cells.Item(138,9).Value = 1338.8
$ws.Cells.Item(138,11).Value = 4016.4
$ws.Cells.Item(138,13).Value = 1123.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2,8).Value = 2955
$ws.Cells.Item(2,9).Value = 3416.25
$ws.Cells.Item(2,11).Value = 3416.25
$ws.Cells.Item(2,13).Value = -3303.25
$ws.Cells.Item(23,8).Value = 23118.75
$ws.Cells.Item(23,10).Value = 24158.334
$ws.Cells.Item(23,12).Value = 24158.334
$ws.Cells.Item(23,14).Value = -24676.334
$ws.Cells.Item(32,8).Value = 7582.0435
$ws.Cells.Item(32,9).Value = 7582.0435
$ws.Cells.Item(32,11).Value = 7582.0435
$ws.Cells.Item(32,13).Value = -7295.0435
$ws.Cells.Item(45,8).Value = 3069.7144
$ws.Cells.Item(45,9).Value = 3069.7144
$ws.Cells.Item(45,10).Value = 0
$ws.Cells.Item(45,11).Value = 3069.7144
$ws.Cells.Item(45,12).Value = 0
$ws.Cells.Item(45,13).ClearContents()
$ws.Cells.Item(45,14).Value = -2692.7144
$ws.Cells.Item(61,8).Value = 1746.4
$ws.Cells.Item(61,9).Value = 1746.4
$ws.Cells.Item(61,11).Value = 1746.4
$ws.Cells.Item(61,13).Value = -1534.4
$ws.Cells.Item(116,8).Value = 2955
$ws.Cells.Item(116,9).Value = 3416.25
$ws.Cells.Item(116,11).Value = 3416.25
$ws.Cells.Item(116,13).Value = -1122.25
$ws.Cells.Item(136,8).Value = 1746.4
$ws.Cells.Item(136,9).Value = 1746.4
$ws.Cells.Item(136,11).Value = 5239.200000000001
$ws.Cells.Item(136,13).Value = -2689.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3,8).Value = 2955
$ws.Cells.Item(3,9).Value = 3416.25
$ws.Cells.Item(3,11).Value = 3416.25
$ws.Cells.Item(3,13).Value = -3302.25
$ws.Cells.Item(116,8).Value = 0
$ws.Cells.Item(116,10).Value = 0
$ws.Cells.Item(116,12).ClearContents()
$ws.Cells.Item(116,14).Value = 0
$ws.Cells.Item(134,8).Value = 2700
$ws.Cells.Item(134,9).Value = 3050
$ws.Cells.Item(134,10).Value = 2000
$ws.Cells.Item(134,11).Value = 9150
$ws.Cells.Item(134,12).Value = 6000
$ws.Cells.Item(134,13).Value = -6615
$ws.Cells.Item(134,14).Value = -11070

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16,8).Value = 1330.2
$ws.Cells.Item(16,10).Value = 1121
$ws.Cells.Item(16,12).Value = 1121
$ws.Cells.Item(16,14).Value = -1695
$ws.Cells.Item(31,8).Value = 2463.6
$ws.Cells.Item(31,9).Value = 2019
$ws.Cells.Item(31,11).Value = 2019
$ws.Cells.Item(31,13).Value = -1724
$ws.Cells.Item(34,8).Value = 2463.6
$ws.Cells.Item(34,9).Value = 2019
$ws.Cells.Item(34,11).Value = 2019
$ws.Cells.Item(34,13).Value = -1817
$ws.Cells.Item(58,8).Value = 1780.4
$ws.Cells.Item(58,9).Value = 1627.375
$ws.Cells.Item(58,10).Value = 2392.5
$ws.Cells.Item(58,11).Value = 1627.375
$ws.Cells.Item(58,12).Value = 2392.5
$ws.Cells.Item(58,13).Value = -1424.375
$ws.Cells.Item(58,14).Value = -2798.5
$ws.Cells.Item(82,8).Value = 16999.5
$ws.Cells.Item(82,9).Value = 14000
$ws.Cells.Item(82,10).Value = 19999
$ws.Cells.Item(82,11).Value = 14000
$ws.Cells.Item(82,12).Value = 19999
$ws.Cells.Item(82,13).Value = -13639
$ws.Cells.Item(82,14).Value = -20721
$ws.Cells.Item(85,8).Value = 16999.5
$ws.Cells.Item(85,9).Value = 14000
$ws.Cells.Item(85,10).Value = 19999
$ws.Cells.Item(85,11).Value = 14000
$ws.Cells.Item(85,12).Value = 19999
$ws.Cells.Item(85,13).Value = -12752
$ws.Cells.Item(85,14).Value = -22495
$ws.Cells.Item(113,8).Value = 1330.2
$ws.Cells.Item(113,10).Value = 1121
$ws.Cells.Item(113,12).Value = 1121
$ws.Cells.Item(113,14).Value = -5461
$ws.Cells.Item(134,8).Value = 7322.4614
$ws.Cells.Item(134,9).Value = 7396.5557
$ws.Cells.Item(134,10).Value = 7155.75
$ws.Cells.Item(134,11).Value = 22189.6671
$ws.Cells.Item(134,12).Value = 21467.25
$ws.Cells.Item(134,13).Value = -19654.6671
$ws.Cells.Item(134,14).Value = -26537.25
$ws.Cells.Item(136,8).Value = 1780.4
$ws.Cells.Item(136,9).Value = 1627.375
$ws.Cells.Item(136,10).Value = 2392.5
$ws.Cells.Item(136,11).Value = 4882.125
$ws.Cells.Item(136,12).Value = 7177.5
$ws.Cells.Item(136,13).Value = -2332.125
$ws.Cells.Item(136,14).Value = -12277.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8,8).Value = 1202488.8
$ws.Cells.Item(8,9).Value = 1202488.8
$ws.Cells.Item(8,11).Value = 3607466.4
$ws.Cells.Item(8,13).Value = -3607327.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102,8).Value = 3744.889
$ws.Cells.Item(102,9).Value = 3430.8572
$ws.Cells.Item(102,10).Value = 4844
$ws.Cells.Item(102,11).Value = 3430.8572
$ws.Cells.Item(102,12).Value = 4844
$ws.Cells.Item(102,13).Value = -1808.8572
$ws.Cells.Item(102,14).Value = -8088
$ws.Cells.Item(126,8).Value = 6455.75
$ws.Cells.Item(126,9).Value = 7303
$ws.Cells.Item(126,10).Value = 3914
$ws.Cells.Item(126,11).Value = 21909
$ws.Cells.Item(126,12).Value = 11742
$ws.Cells.Item(126,13).Value = -19439
$ws.Cells.Item(126,14).Value = -16682
$ws.Cells.Item(132,8).Value = 4222.222
$ws.Cells.Item(132,9).Value = 4000.2856
$ws.Cells.Item(132,11).Value = 12000.8568
$ws.Cells.Item(132,13).Value = -9470.856800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68,8).Value = 2739.2
$ws.Cells.Item(68,10).Value = 2000
$ws.Cells.Item(68,12).Value = 2000
$ws.Cells.Item(68,14).Value = -3498
$ws.Cells.Item(71,8).Value = 2739.2
$ws.Cells.Item(71,10).Value = 2000
$ws.Cells.Item(71,12).Value = 10000
$ws.Cells.Item(71,14).Value = -17488
$ws.Cells.Item(93,8).Value = 506.16666
$ws.Cells.Item(93,10).Value = 392
$ws.Cells.Item(93,12).Value = 392
$ws.Cells.Item(93,14).Value = -2888
$ws.Cells.Item(136,8).Value = 2919.875
$ws.Cells.Item(136,9).Value = 2919.875
$ws.Cells.Item(136,13).Value = -6209.625
